# Adapt column header formatting to respective input file names (#7)
# - rename "<field>_old"  -> "<field>_FV2210"
# - rename "<field>_new"  -> "<field>_FV2304"
# - wrap the sheet's used range in an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21   # A..U
$lastRow = 80   # header + 79 data rows

# 1) Rename the header cells in row 1.
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Text
    if ($header -match "_old$") {
        $cell.Value = ($header -replace "_old$", "_FV2210")
    } elseif ($header -match "_new$") {
        $cell.Value = ($header -replace "_new$", "_FV2304")
    }
}

# 2) Turn the data range into a proper Excel Table, keeping the (renamed)
#    header row as the table's column headers.
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the header row (pane split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
